$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ui")

# New "ui" message rows for the shop recovery items (TournRPG-277).
$newRows = @(
    @{ Row = 26; Text = "回復" },
    @{ Row = 27; Text = "10％回復" },
    @{ Row = 28; Text = "全回復" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $prev = $r - 1

    # Copy the formatting from the row above (A & B columns) so the new
    # rows keep the same styles as the rest of the table.
    $ws.Range("A$prev").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("B$prev").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)

    $ws.Range("A$r").Formula = "=ROW()-2"
    $ws.Range("B$r").Value = $item.Text

    $ws.Rows.Item($r).RowHeight = 20
}
